$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
$cs.Item(1).RGB = 255
